# Update countries & provincias Spain
# Applies updated COVID-19 stats for a new report date, including three
# countries (Bulgaria, Kirguistan, El Salvador) whose alphabetical-ish
# position in the sheet moves up by one row, pushing the following
# country's (unchanged) figures down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - India: simple value update, country unchanged
$ws.Cells.Item(21, 2).Value = 14425
$ws.Cells.Item(21, 3).Value = 73
$ws.Cells.Item(21, 4).Value = 2045
$ws.Cells.Item(21, 5).Value = 11892
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 488

# Row 56 - Tailandia: simple value update, country unchanged
$ws.Cells.Item(56, 2).Value = 2733
$ws.Cells.Item(56, 3).Value = 33
$ws.Cells.Item(56, 4).Value = 1787
$ws.Cells.Item(56, 5).Value = 899
$ws.Cells.Item(56, 6).Value = 61
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 47

# Rows 84-85: Bulgaria moves ahead of Tunez (gets updated figures),
# Tunez shifts down one row keeping its previous figures.
$ws.Cells.Item(84, 1).Value = "Bulgaria"
$ws.Cells.Item(84, 2).Value = 865
$ws.Cells.Item(84, 3).Value = 19
$ws.Cells.Item(84, 4).Value = 153
$ws.Cells.Item(84, 5).Value = 671
$ws.Cells.Item(84, 6).Value = 34
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 41

$ws.Cells.Item(85, 1).Value = "Tunez"
$ws.Cells.Item(85, 2).Value = 864
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 4).Value = 43
$ws.Cells.Item(85, 5).Value = 784
$ws.Cells.Item(85, 6).Value = 89
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 37

# Rows 99-101: Kirguistan moves ahead of Bolivia and Nigeria (gets
# updated figures); Bolivia and Nigeria each shift down one row keeping
# their previous figures.
$ws.Cells.Item(99, 1).Value = "Kirguistan"
$ws.Cells.Item(99, 2).Value = 506
$ws.Cells.Item(99, 3).Value = 17
$ws.Cells.Item(99, 4).Value = 114
$ws.Cells.Item(99, 5).Value = 387
$ws.Cells.Item(99, 6).Value = 5
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 5

$ws.Cells.Item(100, 1).Value = "Bolivia"
$ws.Cells.Item(100, 2).Value = 493
$ws.Cells.Item(100, 3).Value = 28
$ws.Cells.Item(100, 4).Value = 31
$ws.Cells.Item(100, 5).Value = 431
$ws.Cells.Item(100, 6).Value = 3
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 31

$ws.Cells.Item(101, 1).Value = "Nigeria"
$ws.Cells.Item(101, 2).Value = 493
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 159
$ws.Cells.Item(101, 5).Value = 317
$ws.Cells.Item(101, 6).Value = 2
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 17

# Rows 123-124: El Salvador moves ahead of Islas Feroe (gets updated
# figures); Islas Feroe shifts down one row keeping its previous figures.
$ws.Cells.Item(123, 1).Value = "El Salvador"
$ws.Cells.Item(123, 2).Value = 190
$ws.Cells.Item(123, 3).Value = 13
$ws.Cells.Item(123, 4).Value = 43
$ws.Cells.Item(123, 5).Value = 140
$ws.Cells.Item(123, 6).Value = 2
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 7

$ws.Cells.Item(124, 1).Value = "Islas Feroe"
$ws.Cells.Item(124, 2).Value = 184
$ws.Cells.Item(124, 3).Value = 0
$ws.Cells.Item(124, 4).Value = 171
$ws.Cells.Item(124, 5).Value = 13
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 0

# Row 197 - Islas Malvinas: simple value update, country unchanged
$ws.Cells.Item(197, 4).Value = 3
$ws.Cells.Item(197, 5).Value = 8

Write-Host "Applied paises.xlsx updates"
